$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.769.71"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "2.420.15"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'551.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").Value = "'160.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +1.81%  "
$ws.Range("D9").Value = "'0.158"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.75%  "
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("D11").Value = "'4.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("D12").Value = "'0.325"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.24%  "
$ws.Range("D13").Value = "67.733.02"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("D14").Value = "'0.0000168"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.14%  "
$ws.Range("D15").Value = "'22.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("D16").Value = "'10.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.99%  "
$ws.Range("D17").Value = "'335.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.13%  "
$ws.Range("D18").Value = "'6.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.07%  "
$ws.Range("D19").Value = "'3.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").Value = "'1.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.18%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").Value = "'66.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").Value = "'3.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.69%  "
$ws.Range("D24").Value = "'8.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").Value = "0.0₃0810"
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").Value = "'7.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "'421.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.29%  "
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("D30").Value = "'1.59"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.53%  "
$ws.Range("D31").Value = "'161.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.01%  "
$ws.Range("D32").Value = "'18.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").Value = "'17.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").Value = "'0.102"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.99%  "
$ws.Range("D36").Value = "'0.293"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.39%  "
$ws.Range("D37").Value = "'4.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.21%  "
$ws.Range("D38").Value = "'1.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("E39").Value = "  -1.95%  "
$ws.Range("D40").Value = "'2.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.41%  "
$ws.Range("D41").Value = "'3.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("D42").Value = "'128.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.78%  "
$ws.Range("D43").Value = "'0.0709"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("D44").Value = "'0.476"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("E46").Value = "  +0.93%  "
$ws.Range("E47").Value = "  +0.74%  "
$ws.Range("E48").Value = "  -5.43%  "
$ws.Range("D49").Value = "0.0₆0206"
$ws.Range("E49").Value = "  +5.01%  "
$ws.Range("D50").Value = "'16.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.00%  "
$ws.Range("D51").Value = "'4.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.42%  "
